# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) on the player's save-data sheet was being
# populated with a per-game strikeout *count* pulled from an older stat
# ("Strike#"). This regenerates that column with the correct per-game K
# value (s_vals) for every logged outing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> recalculated K value (s_vals)
$sVals = [ordered]@{
    2  = 1
    3  = 0
    4  = 1
    5  = 3
    7  = 0
    8  = 0
    9  = 0
    10 = 2
    11 = 1
    12 = 1
    13 = 2
    14 = 0
    15 = 1
    16 = 1
    17 = 0
    19 = 1
    20 = 1
    21 = 0
    22 = 0
    23 = 3
    24 = 1
    25 = 2
    26 = 1
    27 = 2
    28 = 1
    29 = 6
    30 = 0
    31 = 0
    32 = 1
    33 = 1
    34 = 2
    35 = 0
    37 = 0
    38 = 1
    39 = 2
    40 = 0
    41 = 1
    43 = 1
}

foreach ($row in $sVals.Keys) {
    $ws.Cells.Item($row, 7).Value = $sVals[$row]
}
